$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A70").Value = "State Politics & Policy Quarterly"
$ws.Range("B70").Value = "<a href='https://www.cambridge.org/core/journals/state-politics-and-policy-quarterly'target='_blank'>Short Article</a>"
$ws.Range("C70").Value = "4k words"
$ws.Range("D70").Value = 17

$ws.Range("C71").Select() | Out-Null
